$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert 2 new rows before row 364, shifting existing rows 364:383 down to 366:385.
$ws.Range("A364:R365").Insert()

# Row 364 - Primera
$ws.Cells.Item(364, 1).Value = 8
$ws.Cells.Item(364, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(364, 3).Value = "Coquimbo"
$ws.Cells.Item(364, 4).Value = 44931
$ws.Cells.Item(364, 5).Value = 4
$ws.Cells.Item(364, 6).Value = 100114014
$ws.Cells.Item(364, 7).Value = "Betarraga"
$ws.Cells.Item(364, 8).Value = "Sin especificar"
$ws.Cells.Item(364, 9).Value = "Primera"
$ws.Cells.Item(364, 10).Value = 2000
$ws.Cells.Item(364, 11).Value = 550
$ws.Cells.Item(364, 12).Value = 600
$ws.Cells.Item(364, 13).Value = 575
$ws.Cells.Item(364, 14).Value = "$/paquete 3 unidades"
$ws.Cells.Item(364, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(364, 16).Value = 192
$ws.Cells.Item(364, 17).Value = 3
$ws.Cells.Item(364, 18).Value = "Hortaliza"

# Row 365 - Segunda
$ws.Cells.Item(365, 1).Value = 8
$ws.Cells.Item(365, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(365, 3).Value = "Coquimbo"
$ws.Cells.Item(365, 4).Value = 44931
$ws.Cells.Item(365, 5).Value = 4
$ws.Cells.Item(365, 6).Value = 100114014
$ws.Cells.Item(365, 7).Value = "Betarraga"
$ws.Cells.Item(365, 8).Value = "Sin especificar"
$ws.Cells.Item(365, 9).Value = "Segunda"
$ws.Cells.Item(365, 10).Value = 1560
$ws.Cells.Item(365, 11).Value = 450
$ws.Cells.Item(365, 12).Value = 500
$ws.Cells.Item(365, 13).Value = 475
$ws.Cells.Item(365, 14).Value = "$/paquete 3 unidades"
$ws.Cells.Item(365, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(365, 16).Value = 158
$ws.Cells.Item(365, 17).Value = 3
$ws.Cells.Item(365, 18).Value = "Hortaliza"

# Match date formatting style used by column D elsewhere (s="2")
$ws.Range("D364:D365").NumberFormat = $ws.Range("D366").NumberFormat
